$d = $word.ActiveDocument

# The document currently holds a single paragraph containing just a space.
# Replace its content (and append the new paragraphs after it) by inserting
# raw WordprocessingML over the existing paragraph's range, so that we get
# exactly the paragraph/run shapes the target revision has (including two
# fully-empty <w:p/> separator paragraphs with no run at all).
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$canningText  = "&lt;canningUrl&gt;https://nchfp.uga.edu/how/can_04/peppers.html&lt;/canningUrl&gt;"
$freezingText = "&lt;freezingUrl&gt;https://nchfp.uga.edu/how/freeze/pepper_hot.html&lt;/freezingUrl&gt;"
$dryingText   = "&lt;dryingUrl&gt;https://www.cayennediane.com/how-to-dry-peppers/#:~:text=Cut%20your%20peppers%20in%20half,least%20a%20couple%20of%20inches.&lt;/dryingUrl&gt;"
$pictureText  = "https://i.imgur.com/01iv9ff.jpg"

$xml = "<w:p $wns><w:r><w:t>$canningText</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>$freezingText</w:t></w:r></w:p>" +
       "<w:p $wns><w:r><w:t>$dryingText</w:t></w:r></w:p>" +
       "<w:p $wns/>" +
       "<w:p $wns/>" +
       "<w:p $wns><w:r><w:t>$pictureText</w:t></w:r></w:p>"

$r1.InsertXML($xml)
